$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# --- Header row: shrink italic header font from 11pt (22 half-points) to
# --- 10pt (20 half-points) in all three header cells.
for ($col = 1; $col -le 3; $col++) {
    $cell = $t.Cell(1, $col)
    $r = $cell.Range
    # Exclude the trailing cell-mark character so we don't touch the
    # paragraph mark's run properties.
    $body = $d.Range($r.Start, $r.End - 1)
    $body.Font.Size = 10
}

# --- Data rows: right-align the country-code paragraph in column 1.
$rowCount = $t.Rows.Count
for ($row = 2; $row -le $rowCount; $row++) {
    $cell = $t.Cell($row, 1)
    $cell.Range.Paragraphs.Item(1).Alignment = 2
}
